$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks up front; we will re-add exact set at the end.
$ws.Hyperlinks.Delete()

# Donor cells (outside the used data range) used to re-apply the original
# per-column cell style (format) after writing date-like text, since assigning
# a "dd-mm-yyyy"-looking string via .Value on a General-formatted cell makes
# Excel silently reinterpret it as a real date/time value.
$donorText = $ws.Cells.Item(1, 26)   # Z1  -> style used by columns A,B,C,E,F
$donorNum  = $ws.Cells.Item(1, 27)   # AA1 -> style used by column D
$ws.Cells.Item(2, 2).Copy()
$donorText.PasteSpecial(-4122)
$ws.Cells.Item(2, 4).Copy()
$donorNum.PasteSpecial(-4122)

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $donorText.Copy()
    $cell.PasteSpecial(-4122)
}

# Row 2
Set-TextCell $ws.Cells.Item(2, 1) "10-12-2025"
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(2, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(2, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 3
Set-TextCell $ws.Cells.Item(3, 1) "09-12-2025"
$ws.Cells.Item(3, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(3, 3).Value = "IE07"
$ws.Cells.Item(3, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(3, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(3, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 4
Set-TextCell $ws.Cells.Item(4, 1) "08-12-2025"
$ws.Cells.Item(4, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(4, 3).Value = "IE07"
$ws.Cells.Item(4, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(4, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(4, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 5
Set-TextCell $ws.Cells.Item(5, 1) "07-12-2025"
$ws.Cells.Item(5, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(5, 3).Value = "IE07"
$ws.Cells.Item(5, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(5, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(5, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 6
Set-TextCell $ws.Cells.Item(6, 1) "06-12-2025"
$ws.Cells.Item(6, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(6, 3).Value = "IE07"
$ws.Cells.Item(6, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(6, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(6, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 7
Set-TextCell $ws.Cells.Item(7, 1) "05-12-2025"
$ws.Cells.Item(7, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(7, 3).Value = "IE07"
$ws.Cells.Item(7, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(7, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(7, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 8
Set-TextCell $ws.Cells.Item(8, 1) "04-12-2025"
$ws.Cells.Item(8, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(8, 3).Value = "IE07"
$ws.Cells.Item(8, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(8, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(8, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 9
Set-TextCell $ws.Cells.Item(9, 1) "03-12-2025"
$ws.Cells.Item(9, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(9, 3).Value = "IE07"
$ws.Cells.Item(9, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(9, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(9, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 10
Set-TextCell $ws.Cells.Item(10, 1) "02-12-2025"
$ws.Cells.Item(10, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(10, 3).Value = "IE07"
$ws.Cells.Item(10, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(10, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(10, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 11
Set-TextCell $ws.Cells.Item(11, 1) "01-12-2025"
$ws.Cells.Item(11, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(11, 3).Value = "IE07"
$ws.Cells.Item(11, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(11, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(11, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 12
Set-TextCell $ws.Cells.Item(12, 1) "30-11-2025"
$ws.Cells.Item(12, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(12, 3).Value = "IE07"
$ws.Cells.Item(12, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(12, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(12, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 13
Set-TextCell $ws.Cells.Item(13, 1) "29-11-2025"
$ws.Cells.Item(13, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(13, 3).Value = "IE07"
$ws.Cells.Item(13, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(13, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(13, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 14
Set-TextCell $ws.Cells.Item(14, 1) "28-11-2025"
$ws.Cells.Item(14, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(14, 3).Value = "IE07"
$ws.Cells.Item(14, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(14, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(14, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 15
Set-TextCell $ws.Cells.Item(15, 1) "27-11-2025"
$ws.Cells.Item(15, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(15, 3).Value = "IE07"
$ws.Cells.Item(15, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(15, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(15, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 16
Set-TextCell $ws.Cells.Item(16, 1) "26-11-2025"
$ws.Cells.Item(16, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(16, 3).Value = "IE07"
$ws.Cells.Item(16, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(16, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(16, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 17
Set-TextCell $ws.Cells.Item(17, 1) "25-11-2025"
$ws.Cells.Item(17, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(17, 3).Value = "IE07"
$ws.Cells.Item(17, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(17, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(17, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 18
Set-TextCell $ws.Cells.Item(18, 1) "24-11-2025"
$ws.Cells.Item(18, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(18, 3).Value = "IE07"
$ws.Cells.Item(18, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(18, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(18, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 19
Set-TextCell $ws.Cells.Item(19, 1) "23-11-2025"
$ws.Cells.Item(19, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(19, 3).Value = "IE07"
$ws.Cells.Item(19, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(19, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(19, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 20
Set-TextCell $ws.Cells.Item(20, 1) "22-11-2025"
$ws.Cells.Item(20, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(20, 3).Value = "IE07"
$ws.Cells.Item(20, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(20, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(20, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 21
Set-TextCell $ws.Cells.Item(21, 1) "21-11-2025"
$ws.Cells.Item(21, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(21, 3).Value = "IE07"
$ws.Cells.Item(21, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(21, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(21, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 22
Set-TextCell $ws.Cells.Item(22, 1) "20-11-2025"
$ws.Cells.Item(22, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(22, 3).Value = "IE07"
$ws.Cells.Item(22, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(22, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(22, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 23
Set-TextCell $ws.Cells.Item(23, 1) "19-11-2025"
$ws.Cells.Item(23, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(23, 3).Value = "IE07"
$ws.Cells.Item(23, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(23, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(23, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 24
Set-TextCell $ws.Cells.Item(24, 1) "18-11-2025"
$ws.Cells.Item(24, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(24, 3).Value = "IE07"
$ws.Cells.Item(24, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(24, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(24, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 25
Set-TextCell $ws.Cells.Item(25, 1) "17-11-2025"
$ws.Cells.Item(25, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(25, 3).Value = "IE07"
$ws.Cells.Item(25, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(25, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(25, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 26
Set-TextCell $ws.Cells.Item(26, 1) "16-11-2025"
$ws.Cells.Item(26, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(26, 3).Value = "IE07"
$ws.Cells.Item(26, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(26, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(26, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 27
Set-TextCell $ws.Cells.Item(27, 1) "15-11-2025"
$ws.Cells.Item(27, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(27, 3).Value = "IE07"
$ws.Cells.Item(27, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(27, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(27, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 28
Set-TextCell $ws.Cells.Item(28, 1) "14-11-2025"
$ws.Cells.Item(28, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(28, 3).Value = "IE07"
$ws.Cells.Item(28, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(28, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(28, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 29
Set-TextCell $ws.Cells.Item(29, 1) "13-11-2025"
$ws.Cells.Item(29, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(29, 3).Value = "IE07"
$ws.Cells.Item(29, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(29, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(29, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 30
Set-TextCell $ws.Cells.Item(30, 1) "12-11-2025"
$ws.Cells.Item(30, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(30, 3).Value = "IE07"
$ws.Cells.Item(30, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(30, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(30, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 31
Set-TextCell $ws.Cells.Item(31, 1) "11-11-2025"
$ws.Cells.Item(31, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(31, 3).Value = "IE07"
$ws.Cells.Item(31, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(31, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(31, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 32
Set-TextCell $ws.Cells.Item(32, 1) "10-11-2025"
$ws.Cells.Item(32, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(32, 3).Value = "IE07"
$ws.Cells.Item(32, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(32, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(32, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 33
Set-TextCell $ws.Cells.Item(33, 1) "09-11-2025"
$ws.Cells.Item(33, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(33, 3).Value = "IE07"
$ws.Cells.Item(33, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(33, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(33, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 34
Set-TextCell $ws.Cells.Item(34, 1) "08-11-2025"
$ws.Cells.Item(34, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(34, 3).Value = "IE07"
$ws.Cells.Item(34, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(34, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(34, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 35
Set-TextCell $ws.Cells.Item(35, 1) "07-11-2025"
$ws.Cells.Item(35, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(35, 3).Value = "IE07"
$ws.Cells.Item(35, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(35, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(35, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 36
Set-TextCell $ws.Cells.Item(36, 1) "06-11-2025"
$ws.Cells.Item(36, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(36, 3).Value = "IE07"
$ws.Cells.Item(36, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(36, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(36, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 37
Set-TextCell $ws.Cells.Item(37, 1) "05-11-2025"
$ws.Cells.Item(37, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(37, 3).Value = "IE07"
$ws.Cells.Item(37, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(37, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(37, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 38
Set-TextCell $ws.Cells.Item(38, 1) "04-11-2025"
$ws.Cells.Item(38, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(38, 3).Value = "IE07"
$ws.Cells.Item(38, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(38, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(38, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 39
Set-TextCell $ws.Cells.Item(39, 1) "03-11-2025"
$ws.Cells.Item(39, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(39, 3).Value = "IE07"
$ws.Cells.Item(39, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(39, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(39, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 40
Set-TextCell $ws.Cells.Item(40, 1) "02-11-2025"
$ws.Cells.Item(40, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(40, 3).Value = "IE07"
$ws.Cells.Item(40, 4).Value = 296.05
Set-TextCell $ws.Cells.Item(40, 5) "02-11-2025"
Set-TextCell $ws.Cells.Item(40, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

# Row 41
Set-TextCell $ws.Cells.Item(41, 1) "01-11-2025"
$ws.Cells.Item(41, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(41, 3).Value = "IE07"
$ws.Cells.Item(41, 4).Value = 297.15
Set-TextCell $ws.Cells.Item(41, 5) "01-11-2025"
Set-TextCell $ws.Cells.Item(41, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# Row 42
Set-TextCell $ws.Cells.Item(42, 1) "31-10-2025"
$ws.Cells.Item(42, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(42, 3).Value = "IE07"
$ws.Cells.Item(42, 4).Value = 294.05
Set-TextCell $ws.Cells.Item(42, 5) "30-10-2025"
Set-TextCell $ws.Cells.Item(42, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf"

# Row 43
Set-TextCell $ws.Cells.Item(43, 1) "30-10-2025"
$ws.Cells.Item(43, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(43, 3).Value = "IE07"
$ws.Cells.Item(43, 4).Value = 294.05
Set-TextCell $ws.Cells.Item(43, 5) "30-10-2025"
Set-TextCell $ws.Cells.Item(43, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf"

# Row 44
Set-TextCell $ws.Cells.Item(44, 1) "29-10-2025"
$ws.Cells.Item(44, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(44, 3).Value = "IE07"
$ws.Cells.Item(44, 4).Value = 288.55
Set-TextCell $ws.Cells.Item(44, 5) "25-10-2025"
Set-TextCell $ws.Cells.Item(44, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

# Row 45
Set-TextCell $ws.Cells.Item(45, 1) "28-10-2025"
$ws.Cells.Item(45, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(45, 3).Value = "IE07"
$ws.Cells.Item(45, 4).Value = 288.55
Set-TextCell $ws.Cells.Item(45, 5) "25-10-2025"
Set-TextCell $ws.Cells.Item(45, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

# Row 46
Set-TextCell $ws.Cells.Item(46, 1) "27-10-2025"
$ws.Cells.Item(46, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(46, 3).Value = "IE07"
$ws.Cells.Item(46, 4).Value = 288.55
Set-TextCell $ws.Cells.Item(46, 5) "25-10-2025"
Set-TextCell $ws.Cells.Item(46, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

# Row 47
Set-TextCell $ws.Cells.Item(47, 1) "26-10-2025"
$ws.Cells.Item(47, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(47, 3).Value = "IE07"
$ws.Cells.Item(47, 4).Value = 288.55
Set-TextCell $ws.Cells.Item(47, 5) "25-10-2025"
Set-TextCell $ws.Cells.Item(47, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

# Row 48
Set-TextCell $ws.Cells.Item(48, 1) "25-10-2025"
$ws.Cells.Item(48, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(48, 3).Value = "IE07"
$ws.Cells.Item(48, 4).Value = 288.55
Set-TextCell $ws.Cells.Item(48, 5) "25-10-2025"
Set-TextCell $ws.Cells.Item(48, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

# Row 49
Set-TextCell $ws.Cells.Item(49, 1) "24-10-2025"
$ws.Cells.Item(49, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(49, 3).Value = "IE07"
$ws.Cells.Item(49, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(49, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(49, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 50
Set-TextCell $ws.Cells.Item(50, 1) "23-10-2025"
$ws.Cells.Item(50, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(50, 3).Value = "IE07"
$ws.Cells.Item(50, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(50, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(50, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 51
Set-TextCell $ws.Cells.Item(51, 1) "22-10-2025"
$ws.Cells.Item(51, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(51, 3).Value = "IE07"
$ws.Cells.Item(51, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(51, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(51, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 52
Set-TextCell $ws.Cells.Item(52, 1) "21-10-2025"
$ws.Cells.Item(52, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(52, 3).Value = "IE07"
$ws.Cells.Item(52, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(52, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(52, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 53
Set-TextCell $ws.Cells.Item(53, 1) "20-10-2025"
$ws.Cells.Item(53, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(53, 3).Value = "IE07"
$ws.Cells.Item(53, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(53, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(53, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 54
Set-TextCell $ws.Cells.Item(54, 1) "19-10-2025"
$ws.Cells.Item(54, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(54, 3).Value = "IE07"
$ws.Cells.Item(54, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(54, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(54, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 55
Set-TextCell $ws.Cells.Item(55, 1) "18-10-2025"
$ws.Cells.Item(55, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(55, 3).Value = "IE07"
$ws.Cells.Item(55, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(55, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(55, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 56
Set-TextCell $ws.Cells.Item(56, 1) "17-10-2025"
$ws.Cells.Item(56, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(56, 3).Value = "IE07"
$ws.Cells.Item(56, 4).Value = 282.45
Set-TextCell $ws.Cells.Item(56, 5) "17-10-2025"
Set-TextCell $ws.Cells.Item(56, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

# Row 57
Set-TextCell $ws.Cells.Item(57, 1) "16-10-2025"
$ws.Cells.Item(57, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(57, 3).Value = "IE07"
$ws.Cells.Item(57, 4).Value = 285.05
Set-TextCell $ws.Cells.Item(57, 5) "14-10-2025"
Set-TextCell $ws.Cells.Item(57, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf"

# Row 58
Set-TextCell $ws.Cells.Item(58, 1) "15-10-2025"
$ws.Cells.Item(58, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(58, 3).Value = "IE07"
$ws.Cells.Item(58, 4).Value = 285.05
Set-TextCell $ws.Cells.Item(58, 5) "14-10-2025"
Set-TextCell $ws.Cells.Item(58, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf"

# Row 59
Set-TextCell $ws.Cells.Item(59, 1) "14-10-2025"
$ws.Cells.Item(59, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(59, 3).Value = "IE07"
$ws.Cells.Item(59, 4).Value = 285.05
Set-TextCell $ws.Cells.Item(59, 5) "14-10-2025"
Set-TextCell $ws.Cells.Item(59, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf"

# Row 60
Set-TextCell $ws.Cells.Item(60, 1) "13-10-2025"
$ws.Cells.Item(60, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(60, 3).Value = "IE07"
$ws.Cells.Item(60, 4).Value = 282.85
Set-TextCell $ws.Cells.Item(60, 5) "09-10-2025"
Set-TextCell $ws.Cells.Item(60, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

# Row 61
Set-TextCell $ws.Cells.Item(61, 1) "12-10-2025"
$ws.Cells.Item(61, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(61, 3).Value = "IE07"
$ws.Cells.Item(61, 4).Value = 282.85
Set-TextCell $ws.Cells.Item(61, 5) "09-10-2025"
Set-TextCell $ws.Cells.Item(61, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

# Row 62
Set-TextCell $ws.Cells.Item(62, 1) "11-10-2025"
$ws.Cells.Item(62, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(62, 3).Value = "IE07"
$ws.Cells.Item(62, 4).Value = 282.85
Set-TextCell $ws.Cells.Item(62, 5) "09-10-2025"
Set-TextCell $ws.Cells.Item(62, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

# Row 63
Set-TextCell $ws.Cells.Item(63, 1) "10-10-2025"
$ws.Cells.Item(63, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(63, 3).Value = "IE07"
$ws.Cells.Item(63, 4).Value = 282.85
Set-TextCell $ws.Cells.Item(63, 5) "09-10-2025"
Set-TextCell $ws.Cells.Item(63, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

# Row 64
Set-TextCell $ws.Cells.Item(64, 1) "09-10-2025"
$ws.Cells.Item(64, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(64, 3).Value = "IE07"
$ws.Cells.Item(64, 4).Value = 282.85
Set-TextCell $ws.Cells.Item(64, 5) "09-10-2025"
Set-TextCell $ws.Cells.Item(64, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

# Row 65
Set-TextCell $ws.Cells.Item(65, 1) "08-10-2025"
$ws.Cells.Item(65, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(65, 3).Value = "IE07"
$ws.Cells.Item(65, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(65, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(65, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 66
Set-TextCell $ws.Cells.Item(66, 1) "07-10-2025"
$ws.Cells.Item(66, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(66, 3).Value = "IE07"
$ws.Cells.Item(66, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(66, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(66, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 67
Set-TextCell $ws.Cells.Item(67, 1) "06-10-2025"
$ws.Cells.Item(67, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(67, 3).Value = "IE07"
$ws.Cells.Item(67, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(67, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(67, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 68
Set-TextCell $ws.Cells.Item(68, 1) "05-10-2025"
$ws.Cells.Item(68, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(68, 3).Value = "IE07"
$ws.Cells.Item(68, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(68, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(68, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 69
Set-TextCell $ws.Cells.Item(69, 1) "04-10-2025"
$ws.Cells.Item(69, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(69, 3).Value = "IE07"
$ws.Cells.Item(69, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(69, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(69, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 70
Set-TextCell $ws.Cells.Item(70, 1) "03-10-2025"
$ws.Cells.Item(70, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(70, 3).Value = "IE07"
$ws.Cells.Item(70, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(70, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(70, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 71
Set-TextCell $ws.Cells.Item(71, 1) "02-10-2025"
$ws.Cells.Item(71, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(71, 3).Value = "IE07"
$ws.Cells.Item(71, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(71, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(71, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 72
Set-TextCell $ws.Cells.Item(72, 1) "01-10-2025"
$ws.Cells.Item(72, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(72, 3).Value = "IE07"
$ws.Cells.Item(72, 4).Value = 277.95
Set-TextCell $ws.Cells.Item(72, 5) "01-10-2025"
Set-TextCell $ws.Cells.Item(72, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

# Row 73
Set-TextCell $ws.Cells.Item(73, 1) "30-09-2025"
$ws.Cells.Item(73, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(73, 3).Value = "IE07"
$ws.Cells.Item(73, 4).Value = 274.95
Set-TextCell $ws.Cells.Item(73, 5) "30-09-2025"
Set-TextCell $ws.Cells.Item(73, 6) "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf"

# Row 74
Set-TextCell $ws.Cells.Item(74, 1) "29-09-2025"
$ws.Cells.Item(74, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(74, 3).Value = "IE07"
$ws.Cells.Item(74, 4).Value = 270.25
Set-TextCell $ws.Cells.Item(74, 5) "25-09-2025"
Set-TextCell $ws.Cells.Item(74, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Row 75
Set-TextCell $ws.Cells.Item(75, 1) "28-09-2025"
$ws.Cells.Item(75, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(75, 3).Value = "IE07"
$ws.Cells.Item(75, 4).Value = 270.25
Set-TextCell $ws.Cells.Item(75, 5) "25-09-2025"
Set-TextCell $ws.Cells.Item(75, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Row 76
Set-TextCell $ws.Cells.Item(76, 1) "27-09-2025"
$ws.Cells.Item(76, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(76, 3).Value = "IE07"
$ws.Cells.Item(76, 4).Value = 270.25
Set-TextCell $ws.Cells.Item(76, 5) "25-09-2025"
Set-TextCell $ws.Cells.Item(76, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Row 77
Set-TextCell $ws.Cells.Item(77, 1) "26-09-2025"
$ws.Cells.Item(77, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(77, 3).Value = "IE07"
$ws.Cells.Item(77, 4).Value = 270.25
Set-TextCell $ws.Cells.Item(77, 5) "25-09-2025"
Set-TextCell $ws.Cells.Item(77, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Row 78
Set-TextCell $ws.Cells.Item(78, 1) "25-09-2025"
$ws.Cells.Item(78, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(78, 3).Value = "IE07"
$ws.Cells.Item(78, 4).Value = 270.25
Set-TextCell $ws.Cells.Item(78, 5) "25-09-2025"
Set-TextCell $ws.Cells.Item(78, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

# Row 79
Set-TextCell $ws.Cells.Item(79, 1) "24-09-2025"
$ws.Cells.Item(79, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(79, 3).Value = "IE07"
$ws.Cells.Item(79, 4).Value = 275.25
Set-TextCell $ws.Cells.Item(79, 5) "20-09-2025"
Set-TextCell $ws.Cells.Item(79, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# Row 80
Set-TextCell $ws.Cells.Item(80, 1) "23-09-2025"
$ws.Cells.Item(80, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(80, 3).Value = "IE07"
$ws.Cells.Item(80, 4).Value = 275.25
Set-TextCell $ws.Cells.Item(80, 5) "20-09-2025"
Set-TextCell $ws.Cells.Item(80, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# Row 81
Set-TextCell $ws.Cells.Item(81, 1) "22-09-2025"
$ws.Cells.Item(81, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(81, 3).Value = "IE07"
$ws.Cells.Item(81, 4).Value = 275.25
Set-TextCell $ws.Cells.Item(81, 5) "20-09-2025"
Set-TextCell $ws.Cells.Item(81, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# Row 82
Set-TextCell $ws.Cells.Item(82, 1) "21-09-2025"
$ws.Cells.Item(82, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(82, 3).Value = "IE07"
$ws.Cells.Item(82, 4).Value = 275.25
Set-TextCell $ws.Cells.Item(82, 5) "20-09-2025"
Set-TextCell $ws.Cells.Item(82, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# Row 83
Set-TextCell $ws.Cells.Item(83, 1) "20-09-2025"
$ws.Cells.Item(83, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(83, 3).Value = "IE07"
$ws.Cells.Item(83, 4).Value = 275.25
Set-TextCell $ws.Cells.Item(83, 5) "20-09-2025"
Set-TextCell $ws.Cells.Item(83, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# Row 84
Set-TextCell $ws.Cells.Item(84, 1) "19-09-2025"
$ws.Cells.Item(84, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(84, 3).Value = "IE07"
$ws.Cells.Item(84, 4).Value = 278.95
Set-TextCell $ws.Cells.Item(84, 5) "17-09-2025"
Set-TextCell $ws.Cells.Item(84, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf"

# Row 85
Set-TextCell $ws.Cells.Item(85, 1) "18-09-2025"
$ws.Cells.Item(85, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(85, 3).Value = "IE07"
$ws.Cells.Item(85, 4).Value = 278.95
Set-TextCell $ws.Cells.Item(85, 5) "17-09-2025"
Set-TextCell $ws.Cells.Item(85, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf"

# Row 86
Set-TextCell $ws.Cells.Item(86, 1) "17-09-2025"
$ws.Cells.Item(86, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(86, 3).Value = "IE07"
$ws.Cells.Item(86, 4).Value = 278.95
Set-TextCell $ws.Cells.Item(86, 5) "17-09-2025"
Set-TextCell $ws.Cells.Item(86, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf"

# Row 87
Set-TextCell $ws.Cells.Item(87, 1) "16-09-2025"
$ws.Cells.Item(87, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(87, 3).Value = "IE07"
$ws.Cells.Item(87, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(87, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(87, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 88
Set-TextCell $ws.Cells.Item(88, 1) "15-09-2025"
$ws.Cells.Item(88, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(88, 3).Value = "IE07"
$ws.Cells.Item(88, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(88, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(88, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 89
Set-TextCell $ws.Cells.Item(89, 1) "14-09-2025"
$ws.Cells.Item(89, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(89, 3).Value = "IE07"
$ws.Cells.Item(89, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(89, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(89, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 90
Set-TextCell $ws.Cells.Item(90, 1) "13-09-2025"
$ws.Cells.Item(90, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(90, 3).Value = "IE07"
$ws.Cells.Item(90, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(90, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(90, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 91
Set-TextCell $ws.Cells.Item(91, 1) "12-09-2025"
$ws.Cells.Item(91, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(91, 3).Value = "IE07"
$ws.Cells.Item(91, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(91, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(91, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 92
Set-TextCell $ws.Cells.Item(92, 1) "11-09-2025"
$ws.Cells.Item(92, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(92, 3).Value = "IE07"
$ws.Cells.Item(92, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(92, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(92, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 93
Set-TextCell $ws.Cells.Item(93, 1) "10-09-2025"
$ws.Cells.Item(93, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(93, 3).Value = "IE07"
$ws.Cells.Item(93, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(93, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(93, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 94
Set-TextCell $ws.Cells.Item(94, 1) "09-09-2025"
$ws.Cells.Item(94, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(94, 3).Value = "IE07"
$ws.Cells.Item(94, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(94, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(94, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 95
Set-TextCell $ws.Cells.Item(95, 1) "08-09-2025"
$ws.Cells.Item(95, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(95, 3).Value = "IE07"
$ws.Cells.Item(95, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(95, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(95, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 96
Set-TextCell $ws.Cells.Item(96, 1) "07-09-2025"
$ws.Cells.Item(96, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(96, 3).Value = "IE07"
$ws.Cells.Item(96, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(96, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(96, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 97
Set-TextCell $ws.Cells.Item(97, 1) "06-09-2025"
$ws.Cells.Item(97, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(97, 3).Value = "IE07"
$ws.Cells.Item(97, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(97, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(97, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 98
Set-TextCell $ws.Cells.Item(98, 1) "05-09-2025"
$ws.Cells.Item(98, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(98, 3).Value = "IE07"
$ws.Cells.Item(98, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(98, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(98, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 99
Set-TextCell $ws.Cells.Item(99, 1) "04-09-2025"
$ws.Cells.Item(99, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(99, 3).Value = "IE07"
$ws.Cells.Item(99, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(99, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(99, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 100
Set-TextCell $ws.Cells.Item(100, 1) "03-09-2025"
$ws.Cells.Item(100, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(100, 3).Value = "IE07"
$ws.Cells.Item(100, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(100, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(100, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 101
Set-TextCell $ws.Cells.Item(101, 1) "02-09-2025"
$ws.Cells.Item(101, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(101, 3).Value = "IE07"
$ws.Cells.Item(101, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(101, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(101, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 102
Set-TextCell $ws.Cells.Item(102, 1) "01-09-2025"
$ws.Cells.Item(102, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(102, 3).Value = "IE07"
$ws.Cells.Item(102, 4).Value = 272.05
Set-TextCell $ws.Cells.Item(102, 5) "01-09-2025"
Set-TextCell $ws.Cells.Item(102, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# Row 103
Set-TextCell $ws.Cells.Item(103, 1) "31-08-2025"
$ws.Cells.Item(103, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(103, 3).Value = "IE07"
$ws.Cells.Item(103, 4).Value = 271.05
Set-TextCell $ws.Cells.Item(103, 5) "28-08-2025"
Set-TextCell $ws.Cells.Item(103, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

# Row 104
Set-TextCell $ws.Cells.Item(104, 1) "30-08-2025"
$ws.Cells.Item(104, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(104, 3).Value = "IE07"
$ws.Cells.Item(104, 4).Value = 271.05
Set-TextCell $ws.Cells.Item(104, 5) "28-08-2025"
Set-TextCell $ws.Cells.Item(104, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

# Row 105
Set-TextCell $ws.Cells.Item(105, 1) "29-08-2025"
$ws.Cells.Item(105, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(105, 3).Value = "IE07"
$ws.Cells.Item(105, 4).Value = 271.05
Set-TextCell $ws.Cells.Item(105, 5) "28-08-2025"
Set-TextCell $ws.Cells.Item(105, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

# Row 106
Set-TextCell $ws.Cells.Item(106, 1) "28-08-2025"
$ws.Cells.Item(106, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(106, 3).Value = "IE07"
$ws.Cells.Item(106, 4).Value = 271.05
Set-TextCell $ws.Cells.Item(106, 5) "28-08-2025"
Set-TextCell $ws.Cells.Item(106, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

# Row 107
Set-TextCell $ws.Cells.Item(107, 1) "27-08-2025"
$ws.Cells.Item(107, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(107, 3).Value = "IE07"
$ws.Cells.Item(107, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(107, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(107, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 108
Set-TextCell $ws.Cells.Item(108, 1) "26-08-2025"
$ws.Cells.Item(108, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(108, 3).Value = "IE07"
$ws.Cells.Item(108, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(108, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(108, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 109
Set-TextCell $ws.Cells.Item(109, 1) "25-08-2025"
$ws.Cells.Item(109, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(109, 3).Value = "IE07"
$ws.Cells.Item(109, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(109, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(109, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 110
Set-TextCell $ws.Cells.Item(110, 1) "24-08-2025"
$ws.Cells.Item(110, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(110, 3).Value = "IE07"
$ws.Cells.Item(110, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(110, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(110, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 111
Set-TextCell $ws.Cells.Item(111, 1) "23-08-2025"
$ws.Cells.Item(111, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(111, 3).Value = "IE07"
$ws.Cells.Item(111, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(111, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(111, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 112
Set-TextCell $ws.Cells.Item(112, 1) "22-08-2025"
$ws.Cells.Item(112, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(112, 3).Value = "IE07"
$ws.Cells.Item(112, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(112, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(112, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 113
Set-TextCell $ws.Cells.Item(113, 1) "21-08-2025"
$ws.Cells.Item(113, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(113, 3).Value = "IE07"
$ws.Cells.Item(113, 4).Value = 264.35
Set-TextCell $ws.Cells.Item(113, 5) "21-08-2025"
Set-TextCell $ws.Cells.Item(113, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

# Row 114
Set-TextCell $ws.Cells.Item(114, 1) "20-08-2025"
$ws.Cells.Item(114, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(114, 3).Value = "IE07"
$ws.Cells.Item(114, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(114, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(114, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 115
Set-TextCell $ws.Cells.Item(115, 1) "19-08-2025"
$ws.Cells.Item(115, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(115, 3).Value = "IE07"
$ws.Cells.Item(115, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(115, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(115, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 116
Set-TextCell $ws.Cells.Item(116, 1) "18-08-2025"
$ws.Cells.Item(116, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(116, 3).Value = "IE07"
$ws.Cells.Item(116, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(116, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(116, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 117
Set-TextCell $ws.Cells.Item(117, 1) "17-08-2025"
$ws.Cells.Item(117, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(117, 3).Value = "IE07"
$ws.Cells.Item(117, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(117, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(117, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 118
Set-TextCell $ws.Cells.Item(118, 1) "16-08-2025"
$ws.Cells.Item(118, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(118, 3).Value = "IE07"
$ws.Cells.Item(118, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(118, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(118, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 119
Set-TextCell $ws.Cells.Item(119, 1) "15-08-2025"
$ws.Cells.Item(119, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(119, 3).Value = "IE07"
$ws.Cells.Item(119, 4).Value = 269.45
Set-TextCell $ws.Cells.Item(119, 5) "15-08-2025"
Set-TextCell $ws.Cells.Item(119, 6) "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# Row 120
Set-TextCell $ws.Cells.Item(120, 1) "14-08-2025"
$ws.Cells.Item(120, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(120, 3).Value = "IE07"
$ws.Cells.Item(120, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(120, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(120, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 121
Set-TextCell $ws.Cells.Item(121, 1) "13-08-2025"
$ws.Cells.Item(121, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(121, 3).Value = "IE07"
$ws.Cells.Item(121, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(121, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(121, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 122
Set-TextCell $ws.Cells.Item(122, 1) "12-08-2025"
$ws.Cells.Item(122, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(122, 3).Value = "IE07"
$ws.Cells.Item(122, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(122, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(122, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 123
Set-TextCell $ws.Cells.Item(123, 1) "11-08-2025"
$ws.Cells.Item(123, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(123, 3).Value = "IE07"
$ws.Cells.Item(123, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(123, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(123, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 124
Set-TextCell $ws.Cells.Item(124, 1) "10-08-2025"
$ws.Cells.Item(124, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(124, 3).Value = "IE07"
$ws.Cells.Item(124, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(124, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(124, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 125
Set-TextCell $ws.Cells.Item(125, 1) "09-08-2025"
$ws.Cells.Item(125, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(125, 3).Value = "IE07"
$ws.Cells.Item(125, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(125, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(125, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 126
Set-TextCell $ws.Cells.Item(126, 1) "08-08-2025"
$ws.Cells.Item(126, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(126, 3).Value = "IE07"
$ws.Cells.Item(126, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(126, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(126, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 127
Set-TextCell $ws.Cells.Item(127, 1) "07-08-2025"
$ws.Cells.Item(127, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(127, 3).Value = "IE07"
$ws.Cells.Item(127, 4).Value = 268.25
Set-TextCell $ws.Cells.Item(127, 5) "07-08-2025"
Set-TextCell $ws.Cells.Item(127, 6) "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Row 127 (brand new) column D needs the numeric donor style applied (A,B,C,E,F already get it via Set-TextCell)
$dcell = $ws.Cells.Item(127, 4)
$donorNum.Copy()
$dcell.PasteSpecial(-4122)

# Clean up donor cells so they do not leave stray formatting / expand used range.
$donorText.Clear()
$donorNum.Clear()

# Re-create the hyperlinks for F2:F127, matching each cell's text.
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(21, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(22, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(23, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(24, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(25, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(26, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(27, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(28, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(29, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(30, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(31, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(32, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(33, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(34, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(35, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(36, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(37, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(38, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(39, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(40, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(41, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(42, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(43, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(44, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(45, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(46, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(47, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(48, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(49, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(50, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(51, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(52, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(53, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(54, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(55, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(56, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(57, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(58, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(59, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(60, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(61, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(62, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(63, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(64, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(65, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(66, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(67, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(68, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(69, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(70, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(71, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(72, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(73, 6), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(74, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(75, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(76, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(77, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(78, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(79, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(80, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(81, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(82, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(83, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(84, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(85, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(86, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(87, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(88, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(89, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(90, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(91, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(92, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(93, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(94, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(95, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(96, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(97, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(98, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(99, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(100, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(101, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(102, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(103, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(104, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(105, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(106, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(107, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(108, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(109, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(110, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(111, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(112, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(113, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(114, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(115, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(116, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(117, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(118, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(119, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(120, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(121, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(122, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(123, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(124, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(125, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(126, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(127, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null
